# Insert a new "BEIJING SINODATA TECHNOLOGY CO., LTD." row into the Bank of
# China translation block, pushing the existing rows down by one and
# renumbering the sequential ID column (A) to stay 1..N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 is the first "BANK OF CHINA TIANJIN BRANCH" row today; insert above it
# so the new company slots in right after "ANYANG BRANCH, BANK OF CHINA LIMITED".
$ws.Rows.Item(17).Insert()

$ws.Range("A17").Value = 16.0
$ws.Range("B17").Value = 2.0
$ws.Range("C17").Value = 41370098.0
$ws.Range("D17").Value = "中国银行股份有限公司"
$ws.Range("E17").Value = "BEIJING SINODATA TECHNOLOGY CO., LTD."

# Renumber the sequential "RAP - Translation ID" column so it keeps running
# 1..24 after the insertion (rows 18..25 were rows 17..24 before the insert).
for ($r = 18; $r -le 25; $r++) {
    $ws.Cells.Item($r, 1).Value = [double]($r - 1)
}
